# apiAuto_v10, validate the data in the database using sql
# Populate column F (ActualResponseData) on the "Case2" sheet for rows 2-12,
# mirroring column E (ExpectedResponseData) for most rows — except rows 6 and 12,
# where the actual response differs from the expected one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Case2")

$ws.Range("F2").Value = '{"status":0,"code":"20103","data":null,"msg":"密码不能为空"}'
$ws.Range("F3").Value = '{"status":0,"code":"20103","data":null,"msg":"手机号不能为空"}'
$ws.Range("F4").Value = '{"status":0,"code":"20109","data":null,"msg":"手机号码格式不正确"}'
$ws.Range("F5").Value = '{"status":0,"code":"20108","data":null,"msg":"密码长度必须为6~18"}'
$ws.Range("F6").Value = '{"status":0,"code":"20110","data":null,"msg":"手机号码已被注册"}'
$ws.Range("F7").Value = '{"status":0,"code":"20110","data":null,"msg":"手机号码已被注册"}'
$ws.Range("F8").Value = '{"status":0,"code":"20103","data":null,"msg":"密码不能为空"}'
$ws.Range("F9").Value = '{"status":0,"code":"20103","data":null,"msg":"手机号不能为空"}'
$ws.Range("F10").Value = '{"status":0,"code":"20111","data":null,"msg":"用户名或密码错误"}'
$ws.Range("F11").Value = '{"status":0,"code":"20111","data":null,"msg":"用户名或密码错误"}'
$ws.Range("F12").Value = '{"status":0,"code":"20111","data":null,"msg":"用户名或密码错误"}'
